$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "'305.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.86%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'36.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.044"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.59%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07892"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'1.97%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.239"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'6.37%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'7.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-0.50%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'4.154"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'2.53%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9273"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.72%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.09780"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-0.47%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1871"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'0.51%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.09088"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'4.99%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03751"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'4.44%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09911"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.67%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001430"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-3.66%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005720"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.82%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.463"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.18%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'5.02%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3368"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-2.01%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1320"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-1.01%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'5.082"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'2.68%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'1.78%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04572"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.75%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'-0.36%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004776"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-6.72%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001302"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-7.50%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D39").Value = "'0.01928"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'7.87%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04939"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'5.97%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007824"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.49%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1390"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-0.29%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007815"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'2.31%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002143"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-1.03%"
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'15.18%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006141"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-2.62%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.39%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'52.90%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'-10.36%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'-0.39%"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'-0.39%"
$ws.Range("E51").Style = "Normal"
